$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 2.88
$ws.Range("G3").Value = 4.2
$ws.Range("I3").Value = 2.78
$ws.Range("K3").Value = 5.9
$ws.Range("P3").Value = 1.72
$ws.Range("V3").Value = 1.56
$ws.Range("G4").Value = 2.6
$ws.Range("P4").Value = 1.41
$ws.Range("T5").Value = 2.52
$ws.Range("V5").Value = 1.08
$ws.Range("R6").Value = 1.49
$ws.Range("T6").Value = 1.61
$ws.Range("U6").Value = 2.42
$ws.Range("AE6").Value = 34
$ws.Range("F7").Value = 1.33
$ws.Range("G7").Value = 1.34
$ws.Range("H7").Value = 12
$ws.Range("I7").Value = 12.5
$ws.Range("K7").Value = 6.2
$ws.Range("N7").Value = 5.9
$ws.Range("O7").Value = 1.18
$ws.Range("P7").Value = 2.54
$ws.Range("Q7").Value = 1.56
$ws.Range("R7").Value = 1.66
$ws.Range("S7").Value = 2.44
$ws.Range("T7").Value = 2.06
$ws.Range("U7").Value = 1.83
$ws.Range("X7").Value = 25
$ws.Range("Y7").Value = 42
$ws.Range("Z7").Value = 130
$ws.Range("AB7").Value = 14.5
$ws.Range("AC7").Value = 14.5
$ws.Range("AD7").Value = 46
$ws.Range("AE7").Value = 210
$ws.Range("AG7").Value = 14.5
$ws.Range("AH7").Value = 29
$ws.Range("AI7").Value = 160
$ws.Range("AJ7").Value = 9.800000000000001
$ws.Range("AL7").Value = 55
$ws.Range("AM7").Value = 160
$ws.Range("AN7").Value = 4.4
$ws.Range("AO7").Value = 1000
$ws.Range("I9").Value = 4.7
$ws.Range("J9").Value = 3.65
$ws.Range("V9").Value = 1.27
$ws.Range("AB9").Value = 9.6
$ws.Range("AG9").Value = 1000
$ws.Range("O10").Value = 1.19
$ws.Range("P10").Value = 2.48
$ws.Range("Q10").Value = 1.59
$ws.Range("R10").Value = 1.61
$ws.Range("S10").Value = 2.44
$ws.Range("T10").Value = 2.2
$ws.Range("U10").Value = 1.73
$ws.Range("W10").Value = 4.2
$ws.Range("Y10").Value = 1000
$ws.Range("Z10").Value = 160
$ws.Range("AC10").Value = 17.5
$ws.Range("AD10").Value = 1000
$ws.Range("AE10").Value = 290
$ws.Range("AH10").Value = 1000
$ws.Range("AI10").Value = 210
$ws.Range("AM10").Value = 240
$ws.Range("F11").Value = 1.93
$ws.Range("G11").Value = 1.99
$ws.Range("I11").Value = 4.5
$ws.Range("J11").Value = 3.75
$ws.Range("K11").Value = 3.95
$ws.Range("L11").Value = 1.49
$ws.Range("N11").Value = 3.05
$ws.Range("T11").Value = 1.98
$ws.Range("U11").Value = 1.8
$ws.Range("V11").Value = 1.28
$ws.Range("W11").Value = 2
$ws.Range("Z11").Value = 1000
$ws.Range("AA11").Value = 120
$ws.Range("AB11").Value = 9.800000000000001
$ws.Range("AD11").Value = 20
$ws.Range("AE11").Value = 75
$ws.Range("AI11").Value = 95
$ws.Range("AJ11").Value = 1000
$ws.Range("AK11").Value = 1000
$ws.Range("AN11").Value = 23
$ws.Range("F13").Value = 3.55
$ws.Range("H13").Value = 2.74
$ws.Range("I13").Value = 2.82
$ws.Range("J13").Value = 2.76
$ws.Range("K13").Value = 2.78
$ws.Range("L13").Value = 1.54
$ws.Range("O13").Value = 1.46
$ws.Range("P13").Value = 1.26
$ws.Range("R13").Value = 1.18
$ws.Range("T13").Value = 1.01
$ws.Range("U13").Value = 1.01
$ws.Range("V13").Value = 1.55
$ws.Range("Y13").Value = 10
$ws.Range("Z13").Value = 17
$ws.Range("AA13").Value = 44
$ws.Range("AD13").Value = 14.5
$ws.Range("AE13").Value = 40
$ws.Range("AI13").Value = 70
$ws.Range("AL13").Value = 85
$ws.Range("AM13").Value = 190
$ws.Range("AO13").Value = 44
$ws.Range("H14").Value = 16
$ws.Range("I14").Value = 19.5
$ws.Range("J14").Value = 7.4
$ws.Range("K14").Value = 8
$ws.Range("L14").Value = 1.25
$ws.Range("P14").Value = 2.72
$ws.Range("Q14").Value = 1.51
$ws.Range("R14").Value = 1.69
$ws.Range("S14").Value = 2.26
$ws.Range("U14").Value = 1.68
$ws.Range("W14").Value = 5.3
$ws.Range("Z14").Value = 200
$ws.Range("AC14").Value = 22
$ws.Range("AG14").Value = 1000
$ws.Range("AI14").Value = 250
$ws.Range("AM14").Value = 280
$ws.Range("AN14").Value = 3.75
$ws.Range("F15").Value = 1.86
$ws.Range("K15").Value = 4.2
$ws.Range("L15").Value = 1.44
$ws.Range("N15").Value = 3.6
$ws.Range("O15").Value = 1.34
$ws.Range("P15").Value = 1.96
$ws.Range("Q15").Value = 1.95
$ws.Range("R15").Value = 1.33
$ws.Range("S15").Value = 3.35
$ws.Range("T15").Value = 1.83
$ws.Range("U15").Value = 1.96
$ws.Range("AA15").Value = 130
$ws.Range("AB15").Value = 11
$ws.Range("AM15").Value = 140
$ws.Range("F16").Value = 2.22
$ws.Range("G16").Value = 2.26
$ws.Range("H16").Value = 3.45
$ws.Range("I16").Value = 3.65
$ws.Range("L16").Value = 1.43
$ws.Range("O16").Value = 1.32
$ws.Range("V16").Value = 1.38
$ws.Range("W16").Value = 1.79
$ws.Range("AJ16").Value = 29
$ws.Range("AK16").Value = 25
$ws.Range("AO16").Value = 44
$ws.Range("H17").Value = 8.199999999999999
$ws.Range("I17").Value = 9.6
$ws.Range("J17").Value = 4.5
$ws.Range("K17").Value = 5
$ws.Range("P17").Value = 1.98
$ws.Range("T17").Value = 2.08
$ws.Range("U17").Value = 1.76
$ws.Range("X17").Value = 18.5
$ws.Range("AD17").Value = 36
$ws.Range("AG17").Value = 11
$ws.Range("AI17").Value = 170
$ws.Range("F18").Value = 1.18
$ws.Range("I18").Value = 28
$ws.Range("K18").Value = 9.199999999999999
$ws.Range("P18").Value = 2.76
$ws.Range("AM18").Value = 350
$ws.Range("AN18").Value = 3.55
$ws.Range("J19").Value = 3.65
$ws.Range("X19").Value = 16.5
$ws.Range("AB19").Value = 16
$ws.Range("AL19").Value = 46
$ws.Range("G20").Value = 1.32
$ws.Range("H20").Value = 12
$ws.Range("J20").Value = 6.4
$ws.Range("R20").Value = 1.78
$ws.Range("U20").Value = 2.06
$ws.Range("AC20").Value = 15.5
$ws.Range("AE20").Value = 180
$ws.Range("AN20").Value = 3.9
$ws.Range("H21").Value = 6.2
$ws.Range("I21").Value = 6.8
$ws.Range("K21").Value = 4.8
$ws.Range("V21").Value = 1.17
$ws.Range("G22").Value = 1.26
$ws.Range("K22").Value = 8.800000000000001
$ws.Range("P22").Value = 3.35
$ws.Range("R22").Value = 1.99
$ws.Range("AC22").Value = 21
$ws.Range("AI22").Value = 150
$ws.Range("F23").Value = 1.63
$ws.Range("G23").Value = 1.71
$ws.Range("I23").Value = 5.9
$ws.Range("J23").Value = 4.3
$ws.Range("K23").Value = 4.8
$ws.Range("M23").Value = 1.04
$ws.Range("N23").Value = 5
$ws.Range("O23").Value = 1.22
$ws.Range("P23").Value = 2.38
$ws.Range("Q23").Value = 1.65
$ws.Range("R23").Value = 1.54
$ws.Range("S23").Value = 2.66
$ws.Range("T23").Value = 1.71
$ws.Range("U23").Value = 2.22
$ws.Range("W23").Value = 2.42
$ws.Range("Y23").Value = 27
$ws.Range("AF23").Value = 13
$ws.Range("AK23").Value = 18.5
$ws.Range("AN23").Value = 8.6
